$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K9").Value = 1.58
$ws.Range("J10").Value = 1.57
$ws.Range("K10").Value = 2.35
$ws.Range("AG10").Value = 1.17
$ws.Range("AH10").Value = 5
$ws.Range("G12").Value = 2.3
$ws.Range("I12").Value = 3.05
$ws.Range("J12").Value = 2.35
$ws.Range("O12").Value = 9.75
$ws.Range("Q12").Value = 23
$ws.Range("R12").Value = 23
$ws.Range("Y12").Value = 7.2
$ws.Range("AA12").Value = 11.75
$ws.Range("AB12").Value = 37
$ws.Range("AC12").Value = 32
$ws.Range("AD12").Value = 50
$ws.Range("AJ14").Value = 1.63
$ws.Range("G15").Value = 1.26
$ws.Range("H15").Value = 4.5
$ws.Range("I15").Value = 10.75
$ws.Range("J15").Value = 1.8
$ws.Range("K15").Value = 1.8
$ws.Range("L15").Value = 1.37
$ws.Range("M15").Value = 2.52
$ws.Range("N15").Value = 4.7
$ws.Range("O15").Value = 4.45
$ws.Range("Q15").Value = 5.9
$ws.Range("R15").Value = 10.25
$ws.Range("S15").Value = 32
$ws.Range("T15").Value = 9.5
$ws.Range("U15").Value = 8.25
$ws.Range("V15").Value = 24
$ws.Range("W15").Value = 120
$ws.Range("Y15").Value = 18.5
$ws.Range("Z15").Value = 60
$ws.Range("AA15").Value = 28
$ws.Range("AD15").Value = 110
$ws.Range("G16").Value = 2.2
$ws.Range("H16").Value = 3.35
$ws.Range("J16").Value = 1.83
$ws.Range("K16").Value = 1.78
$ws.Range("L16").Value = 1.37
$ws.Range("M16").Value = 2.52
$ws.Range("N16").Value = 6.9
$ws.Range("O16").Value = 9.25
$ws.Range("P16").Value = 7.7
$ws.Range("Q16").Value = 17.5
$ws.Range("R16").Value = 14.5
$ws.Range("S16").Value = 22
$ws.Range("T16").Value = 10
$ws.Range("U16").Value = 5.7
$ws.Range("V16").Value = 11.5
$ws.Range("W16").Value = 45
$ws.Range("X16").Value = 300
$ws.Range("Y16").Value = 7.8
$ws.Range("Z16").Value = 12
$ws.Range("AC16").Value = 19
$ws.Range("AD16").Value = 25
$ws.Range("G17").Value = 2.75
$ws.Range("I17").Value = 2.63
$ws.Range("AC17").Value = 23
$ws.Range("AI17").Value = 1.87
$ws.Range("AJ17").Value = 1.87
$ws.Range("J18").Value = 1.98
$ws.Range("K18").Value = 1.83
$ws.Range("AE18").Value = 1.06
$ws.Range("AF18").Value = 10
$ws.Range("AG18").Value = 1.29
$ws.Range("AH18").Value = 3.5
$ws.Range("AI18").Value = 1.67
$ws.Range("K19").Value = 1.67
$ws.Range("K20").Value = 1.58
$ws.Range("AG21").Value = 1.38
$ws.Range("AH21").Value = 2.6
$ws.Range("AE22").Value = 1.06
$ws.Range("AG22").Value = 1.29
$ws.Range("J23").Value = 2.6
$ws.Range("K23").Value = 1.44
$ws.Range("AE23").Value = 1.11
$ws.Range("AH23").Value = 2.5
$ws.Range("G24").Value = 1.62

Write-Host "Applied all odds updates"
